$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values for B2:K21 (row => comma-separated values for columns B..K)
$rowData = @{
  2 = @(-19.15386123515478,-0.9990381826659661,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  3 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  4 = @(-19.15386123515478,-0.858679734635326,-0.30616460333648,-19.15386123515478,3.994246166533539,-19.15386123515478,2.39893214252466,-19.15386123515478,3.093679537043958,-19.15386123515478)
  5 = @(-19.15386123515478,0.02313241212271987,-19.15386123515478,-19.15386123515478,-19.15386123515478,3.595498589364069,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  6 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  7 = @(3.306733892762075,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  8 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,1.065417202273131,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  9 = @(3.336959534865764,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  10 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,2.14402731005026)
  11 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,2.331373056875461,-19.15386123515478,0.9230339354371069,-19.15386123515478,-19.15386123515478,-19.15386123515478,1.291098022481198)
  12 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  13 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,1.998768358226376,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,0.1544889059976452,3.080301589690726)
  14 = @(-19.15386123515478,-19.15386123515478,0.6136006952612457,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,1.387879161740684)
  15 = @(-19.15386123515478,-19.15386123515478,-0.8262164958982299,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  16 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,1.775860223923641,-19.15386123515478)
  17 = @(-19.15386123515478,-0.1889760972638846,-0.8006584001494811,-19.15386123515478,-19.15386123515478,-19.15386123515478,1.895788595100605,-19.15386123515478,1.528177494665583,-19.15386123515478)
  18 = @(-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,-19.15386123515478,2.315898315293611,-19.15386123515478,2.014928998994455,-19.15386123515478)
  19 = @(-19.15386123515478,-19.15386123515478,3.007480119206345,-19.15386123515478,-19.15386123515478,-19.15386123515478,1.613939312153184,-19.15386123515478,-19.15386123515478,-19.15386123515478)
  20 = @(-19.15386123515478,3.403511077879447,3.084308185378482,-19.15386123515478,2.022777610016041,-19.15386123515478,0.4054654372468621,4.32192574519106,-19.15386123515478,1.041271833032086)
  21 = @(-19.15386123515478,2.69443657143957,-19.15386123515478,3.150197305743947,-19.15386123515478,2.588784138606857,0.7137568633776433,-19.15386123515478,-19.15386123515478,-19.15386123515478)
}

foreach ($r in $rowData.Keys) {
  $vals = $rowData[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
  }
}
